$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.762195110321045
$ws.Range("B1").Value = 2.464328527450562
$ws.Range("C1").Value = 1.74525785446167
$ws.Range("D1").Value = 1.582785844802856
$ws.Range("E1").Value = 1.607695579528809
